$d = $word.ActiveDocument

function Find-ParagraphContaining($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# Split the run(s) covering $doc character offset $pos into two runs at
# that exact point, without touching formatting: add + immediately delete a
# bookmark there. This mirrors how the source document ends up with two
# adjacent runs (same rPr) instead of one run after an edit that otherwise
# gets silently coalesced.
$script:splitCounter = 0
function Split-RunAt($doc, [int]$pos) {
    $r = $doc.Range($pos, $pos)
    $script:splitCounter = $script:splitCounter + 1
    $bmName = "tmpsplit" + $script:splitCounter
    $doc.Bookmarks.Add($bmName, $r) | Out-Null
    $doc.Bookmarks($bmName).Delete()
}

# -----------------------------------------------------------------
# Change 1: "Note:" italic sentence - drop the "// more needed " tail.
#   before (3 runs): "required to changes component values. //"
#                   + "more needed" + " "
#   after  (2 runs): "requir" + "ed to changes component values."
# -----------------------------------------------------------------
$p1 = Find-ParagraphContaining $d "required to changes component values."

$full1 = $p1.Range.Text
$tail = "ed to changes component values. //more needed "
$idxTail = $full1.IndexOf($tail)
$startTail = $p1.Range.Start + $idxTail
$endTail = $startTail + $tail.Length
$rTail = $d.Range($startTail, $endTail)
$rTail.Text = "ed to changes component values."

# Re-create the run boundary right after "requir" so it keeps living in its
# own run instead of being silently coalesced with the text that follows.
$p1Start = $p1.Range.Start
$splitPos1 = $p1Start + $idxTail
Split-RunAt $d $splitPos1

# -----------------------------------------------------------------
# Change 2: drop the hard-coded postgres/S0laT0nga credentials and point
# the reader at their IT administrator instead.
#   before (multiple runs): "You will be asked for a username and password"
#       + ", if you store ... input username and password. If prompted,
#          input username as " + "'" + "postgres" + "'" + " and password "
#       + "'S0laT0nga'" + "."
#   after (3 runs): "You will be asked for a username and password"
#       + ", if you store ... input user"
#       + "name and password. If prompted please ask your IT administrator
#          for login details."
# -----------------------------------------------------------------
$p2 = Find-ParagraphContaining $d "You will be asked for a username and password"

$full2 = $p2.Range.Text
$target2 = "username and password. If prompted, input username as 'postgres' and password 'S0laT0nga'."
$idx2 = $full2.IndexOf($target2)
$p2Start = $p2.Range.Start
$start2 = $p2Start + $idx2
$end2 = $start2 + $target2.Length
$r2 = $d.Range($start2, $end2)
$r2.Text = "username and password. If prompted please ask your IT administrator for login details."

# Restore the run boundary right after "...and password" (end of the
# original, untouched first run) so it does not get coalesced with the new
# replacement text.
$boundary1 = "You will be asked for a username and password"
$boundary1Len = $boundary1.Length
$splitPos2 = $p2Start + $boundary1Len
Split-RunAt $d $splitPos2

# Restore a second run boundary right after "...input user" / before
# "name and password..." to match the two-run split used for the
# replacement sentence.
$full2b = $p2.Range.Text
$suffix2 = "name and password. If prompted please ask your IT administrator for login details."
$idxSuffix2 = $full2b.IndexOf($suffix2)
$splitPos3 = $p2Start + $idxSuffix2
Split-RunAt $d $splitPos3

Write-Output "Change 1 paragraph: $($p1.Range.Text)"
Write-Output "Change 2 paragraph: $($p2.Range.Text)"
